$d = $word.ActiveDocument
$rng = $d.Content
$found = $rng.Find.Execute("Achievement Highlights:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Output "WordOpenXML:"
Write-Output $rng.WordOpenXML
